# Automatische test-sync: 2025-07-31 21:27:50
# Adds a new log row (row 5) to the "Logs" sheet for Testmail #3, extends the
# conditional formatting ranges that covered rows 2-4 to also cover row 5,
# and bumps the "Overig" count on the "Dashboard" sheet from 2 to 3.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row with the new test-mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Kun jij dit afhandelen?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #3: Kun jij dit afhandelen?"
$logs.Range("D5").Value = "Overig"
$logs.Range("E5").Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$logs.Range("F5").Value = "2025-07-31 21:27:34"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Ja"
$logs.Range("I5").Value = "Nee"
$logs.Range("J5").Value = "Nee"

# --- Extend the existing conditional-formatting blocks to include row 5 ---
$logs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))
$logs.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J5"))

# --- Dashboard sheet: bump the "Overig" tally to reflect the new row ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 3
